$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value2 = 23709; $ws.Range("D2").Value2 = 34610773
$ws.Range("C3").Value2 = 59767; $ws.Range("D3").Value2 = 88466880
$ws.Range("C4").Value2 = 20282; $ws.Range("D4").Value2 = 30260617
$ws.Range("C5").Value2 = 5359; $ws.Range("D5").Value2 = 8016370
$ws.Range("C6").Value2 = 1064; $ws.Range("D6").Value2 = 1593697
$ws.Range("C10").Value2 = 25399; $ws.Range("D10").Value2 = 35093436
$ws.Range("C11").Value2 = 6225; $ws.Range("D11").Value2 = 9107750
$ws.Range("C12").Value2 = 17560; $ws.Range("D12").Value2 = 25969702
$ws.Range("C13").Value2 = 5467; $ws.Range("D13").Value2 = 8162839
$ws.Range("C14").Value2 = 1293; $ws.Range("D14").Value2 = 1934129
$ws.Range("C15").Value2 = 234; $ws.Range("D15").Value2 = 348266
$ws.Range("C17").Value2 = 6378; $ws.Range("D17").Value2 = 8641015
$ws.Range("C18").Value2 = 8679; $ws.Range("D18").Value2 = 12651874
$ws.Range("C19").Value2 = 21491; $ws.Range("D19").Value2 = 31808739
$ws.Range("C20").Value2 = 6821; $ws.Range("D20").Value2 = 10193528
$ws.Range("C21").Value2 = 1640; $ws.Range("D21").Value2 = 2454802
$ws.Range("C24").Value2 = 7421; $ws.Range("D24").Value2 = 10145428
$ws.Range("C25").Value2 = 4910; $ws.Range("D25").Value2 = 7162134
$ws.Range("C26").Value2 = 15145; $ws.Range("D26").Value2 = 22406007
$ws.Range("C27").Value2 = 5153; $ws.Range("D27").Value2 = 7703948
$ws.Range("C28").Value2 = 1240; $ws.Range("D28").Value2 = 1859491
$ws.Range("C31").Value2 = 5242; $ws.Range("D31").Value2 = 7042118
$ws.Range("C32").Value2 = 1736; $ws.Range("D32").Value2 = 2503380
$ws.Range("C33").Value2 = 4571; $ws.Range("D33").Value2 = 6716126
$ws.Range("C34").Value2 = 1855; $ws.Range("D34").Value2 = 2761091
$ws.Range("C35").Value2 = 479; $ws.Range("D35").Value2 = 715041
$ws.Range("C38").Value2 = 1172; $ws.Range("D38").Value2 = 1594448
$ws.Range("C39").Value2 = 11040; $ws.Range("D39").Value2 = 16101034
$ws.Range("C40").Value2 = 34032; $ws.Range("D40").Value2 = 50320712
$ws.Range("C41").Value2 = 12606; $ws.Range("D41").Value2 = 18834742
$ws.Range("C42").Value2 = 3483; $ws.Range("D42").Value2 = 5214266
$ws.Range("C43").Value2 = 604; $ws.Range("D43").Value2 = 904936
$ws.Range("C46").Value2 = 10374; $ws.Range("D46").Value2 = 14129822
$ws.Range("C47").Value2 = 995; $ws.Range("D47").Value2 = 1440597
$ws.Range("C48").Value2 = 3737; $ws.Range("D48").Value2 = 5511619
$ws.Range("C49").Value2 = 1396; $ws.Range("D49").Value2 = 2086964
$ws.Range("C50").Value2 = 431; $ws.Range("D50").Value2 = 644000
$ws.Range("C52").Value2 = 2443; $ws.Range("D52").Value2 = 3400256
$ws.Range("C54").Value2 = 953; $ws.Range("D54").Value2 = 1411492
$ws.Range("C55").Value2 = 385; $ws.Range("D55").Value2 = 575476
$ws.Range("C56").Value2 = 132; $ws.Range("D56").Value2 = 197878
$ws.Range("C58").Value2 = 462; $ws.Range("D58").Value2 = 656481
$ws.Range("C59").Value2 = 10018; $ws.Range("D59").Value2 = 14555113
$ws.Range("C60").Value2 = 30510; $ws.Range("D60").Value2 = 45018621
$ws.Range("C61").Value2 = 10589; $ws.Range("D61").Value2 = 15826158
$ws.Range("C62").Value2 = 2932; $ws.Range("D62").Value2 = 4388138
$ws.Range("C63").Value2 = 519; $ws.Range("D63").Value2 = 778139
$ws.Range("C66").Value2 = 9911; $ws.Range("D66").Value2 = 13283986
$ws.Range("C67").Value2 = 2720; $ws.Range("D67").Value2 = 3971418
$ws.Range("C68").Value2 = 7399; $ws.Range("D68").Value2 = 10920685
$ws.Range("C69").Value2 = 2623; $ws.Range("D69").Value2 = 3918983
$ws.Range("C70").Value2 = 861; $ws.Range("D70").Value2 = 1289549
$ws.Range("C73").Value2 = 2864; $ws.Range("D73").Value2 = 3902465
$ws.Range("C75").Value2 = 2998; $ws.Range("D75").Value2 = 4431546
$ws.Range("C76").Value2 = 1194; $ws.Range("D76").Value2 = 1787159
$ws.Range("C80").Value2 = 1772; $ws.Range("D80").Value2 = 2379211
$ws.Range("C81").Value2 = 308; $ws.Range("D81").Value2 = 458189
$ws.Range("C82").Value2 = 106; $ws.Range("D82").Value2 = 158610
$ws.Range("C86").Value2 = 7036; $ws.Range("D86").Value2 = 10292210
$ws.Range("C87").Value2 = 20222; $ws.Range("D87").Value2 = 29918100
$ws.Range("C88").Value2 = 6639; $ws.Range("D88").Value2 = 9923215
$ws.Range("C89").Value2 = 1758; $ws.Range("D89").Value2 = 2632655
$ws.Range("C93").Value2 = 6312; $ws.Range("D93").Value2 = 8513422
$ws.Range("C94").Value2 = 19336; $ws.Range("D94").Value2 = 28080992
$ws.Range("C95").Value2 = 44886; $ws.Range("D95").Value2 = 66214600
$ws.Range("C96").Value2 = 14352; $ws.Range("D96").Value2 = 21429010
$ws.Range("C97").Value2 = 3837; $ws.Range("D97").Value2 = 5741384
$ws.Range("C98").Value2 = 655; $ws.Range("D98").Value2 = 980862
$ws.Range("C101").Value2 = 16532; $ws.Range("D101").Value2 = 22461718
$ws.Range("C102").Value2 = 22106; $ws.Range("D102").Value2 = 32145448
$ws.Range("C103").Value2 = 50014; $ws.Range("D103").Value2 = 73696602
$ws.Range("C104").Value2 = 15625; $ws.Range("D104").Value2 = 23302382
$ws.Range("C105").Value2 = 4002; $ws.Range("D105").Value2 = 5979314
$ws.Range("C106").Value2 = 647; $ws.Range("D106").Value2 = 967554
$ws.Range("C109").Value2 = 19590; $ws.Range("D109").Value2 = 26434297
$ws.Range("C110").Value2 = 8604; $ws.Range("D110").Value2 = 12568021
$ws.Range("C111").Value2 = 22279; $ws.Range("D111").Value2 = 32978078
$ws.Range("C112").Value2 = 7728; $ws.Range("D112").Value2 = 11537814
$ws.Range("C113").Value2 = 1871; $ws.Range("D113").Value2 = 2799094
$ws.Range("C114").Value2 = 268; $ws.Range("D114").Value2 = 399262
$ws.Range("C117").Value2 = 7013; $ws.Range("D117").Value2 = 9549983
$ws.Range("C118").Value2 = 21522; $ws.Range("D118").Value2 = 31296750
$ws.Range("C119").Value2 = 53073; $ws.Range("D119").Value2 = 78280449
$ws.Range("C120").Value2 = 15977; $ws.Range("D120").Value2 = 23861272
$ws.Range("C121").Value2 = 3971; $ws.Range("D121").Value2 = 5939997
$ws.Range("C122").Value2 = 810; $ws.Range("D122").Value2 = 1213212
$ws.Range("C125").Value2 = 18303; $ws.Range("D125").Value2 = 25159506
$ws.Range("C126").Value2 = 29528; $ws.Range("D126").Value2 = 43250772
$ws.Range("C127").Value2 = 88846; $ws.Range("D127").Value2 = 131609166
$ws.Range("C128").Value2 = 39424; $ws.Range("D128").Value2 = 58926655
$ws.Range("C129").Value2 = 12284; $ws.Range("D129").Value2 = 18399251
$ws.Range("C130").Value2 = 2497; $ws.Range("D130").Value2 = 3739909
$ws.Range("C134").Value2 = 29045; $ws.Range("D134").Value2 = 40508905
